# Refresh cryptocurrency Price / Volume(1h) figures scraped into "Sheet1"
# (scheduled GitHub Actions data update); also corrects the row 44/45
# coin ordering (PolygonEcosystemToken now ranks above USDe).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.228.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.54%  "
$ws.Range("D3").Value = "'3.110.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'219.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.64%  "
$ws.Range("D6").Value = "'622.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("D7").Value = "'0.379"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.71%  "
$ws.Range("D8").Value = "'0.970"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +21.04%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "'3.108.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").Value = "'0.721"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +21.23%  "
$ws.Range("E12").Value = "  +5.18%  "
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.95%  "
$ws.Range("D14").Value = "'34.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.20%  "
$ws.Range("D15").Value = "'91.132.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "'3.691.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").Value = "'3.128.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("D19").Value = "'3.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +13.54%  "
$ws.Range("E20").Value = "  +10.19%  "
$ws.Range("D21").Value = "'14.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.21%  "
$ws.Range("D22").Value = "'435.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.69%  "
$ws.Range("D23").Value = "'8.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +6.29%  "
$ws.Range("D25").Value = "'6.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.93%  "
$ws.Range("D26").Value = "'87.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.69%  "
$ws.Range("D27").Value = "'12.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.75%  "
$ws.Range("D28").Value = "'3.285.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("D31").Value = "'9.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.47%  "
$ws.Range("D32").Value = "'525.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.80%  "
$ws.Range("D33").Value = "'0.888"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -18.14%  "
$ws.Range("D34").Value = "'3.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.45%  "
$ws.Range("D35").Value = "'7.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.20%  "
$ws.Range("E36").Value = "  +9.43%  "
$ws.Range("D37").Value = "'23.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.39%  "
$ws.Range("E38").Value = "  +3.76%  "
$ws.Range("E39").Value = "  +3.52%  "
$ws.Range("D40").Value = "'0.0868"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +25.73%  "
$ws.Range("D41").Value = "'22.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").Value = "'0.151"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.06%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "'0.385"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.21%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'1.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.45%  "
$ws.Range("D47").Value = "'146.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("D48").Value = "'43.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("E49").Value = "  +9.70%  "
$ws.Range("D50").Value = "'166.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.08%  "
$ws.Range("E51").Value = "  +6.81%  "
